# Refresh the "想去人数" (want-to-go headcount) figures in column F
# across the sheets that carry this table, matching the regenerated scrape output.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 625
$ws.Range("F7").Value = 3164
$ws.Range("F9").Value = 495
$ws.Range("F10").Value = 2036
$ws.Range("F11").Value = 438
$ws.Range("F12").Value = 372
$ws.Range("F15").Value = 208
$ws.Range("F16").Value = 1021
$ws.Range("F20").Value = 3862
$ws.Range("F21").Value = 1203
$ws.Range("F22").Value = 3054
$ws.Range("F24").Value = 50
$ws.Range("F25").Value = 2599
$ws.Range("F26").Value = 4420
$ws.Range("F30").Value = 2956
$ws.Range("F33").Value = 66
$ws.Range("F34").Value = 52
$ws.Range("F36").Value = 1069
$ws.Range("F37").Value = 1312
$ws.Range("F39").Value = 1166
$ws.Range("F40").Value = 753
$ws.Range("F42").Value = 689
$ws.Range("F45").Value = 158
$ws.Range("F49").Value = 3640

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 649

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 625
$ws.Range("F7").Value = 3164
$ws.Range("F8").Value = 495
$ws.Range("F10").Value = 2036
$ws.Range("F11").Value = 438
$ws.Range("F12").Value = 372
$ws.Range("F18").Value = 208
$ws.Range("F19").Value = 1021
$ws.Range("F23").Value = 3862
$ws.Range("F25").Value = 1203
$ws.Range("F27").Value = 3054
$ws.Range("F28").Value = 2599
$ws.Range("F29").Value = 4420
$ws.Range("F32").Value = 2956
$ws.Range("F34").Value = 1069
$ws.Range("F35").Value = 1312
$ws.Range("F37").Value = 1166
$ws.Range("F38").Value = 753
$ws.Range("F44").Value = 158
$ws.Range("F48").Value = 3640
